$d = $word.ActiveDocument

# 1. Remove the red (FF0000) font color from the "Remake H-bridge symbol..."
#    paragraph, leaving the bold formatting intact. Iterating via the
#    Paragraph's Range (rather than a plain Range from Find) also picks up
#    the paragraph-mark run properties (pPr/rPr), matching the paragraph
#    mark color removal seen in the target.
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "Remake H-bridge symbol*") {
        $p.Range.Font.Color = -16777216
    }
}

# 2. Move the "_GoBack" bookmark from the end of that paragraph to inside
#    the word "Tufte" (between "Tuf" and "te") in the Edward Tufte
#    paragraph - this is what Word does automatically to mark the last
#    edit location, and it also removes the old occurrence of the
#    bookmark since bookmark names are unique per document.
$rng = $d.Content
$found = $rng.Find.Execute("Edward Tuf", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$target = $d.Range($rng.End, $rng.End)
$d.Bookmarks.Add("_GoBack", $target)
